$wb = $excel.ActiveWorkbook

# --- Set up selection on the existing "JAVA SOCKETS RESULTS" sheet -------
$wsJava = $wb.Worksheets("JAVA SOCKETS RESULTS")
$wsJava.Range("A3:A13").Select()

# --- Add the new "MEM MAP RESULTS" sheet at the end of the workbook ------
$wsMem = $wb.Worksheets.Add($null, $wb.Worksheets($wb.Worksheets.Count))
$wsMem.Name = "MEM MAP RESULTS"

# Title row (merged A1:E1), reusing the same header style/text as the other
# "results" sheets.
$wsMem.Range("A1").Value = "Java Sockets (No JNI)"
$wsMem.Range("A1:E1").Merge()

# Column headers
$wsMem.Range("B2").Value = "40 BYTES "
$wsMem.Range("C2").Value = "400 BYTES"
$wsMem.Range("D2").Value = "4000 BYTES"
$wsMem.Range("E2").Value = "40 000 BYTES"

# Data rows 3-12
$data = @(
    @(1, 106100, 106855, 103457, 81557),
    @(2, 89110, 99682, 76649, 118561),
    @(3, 83823, 86089, 74761, 71363),
    @(4, 71741, 74383, 350395, 131021),
    @(5, 69474, 76271, 95150, 74383),
    @(6, 79292, 78537, 103079, 79291),
    @(7, 118560, 152920, 103080, 72873),
    @(8, 75516, 126867, 118938, 76649),
    @(9, 79292, 72117, 86844, 79292),
    @(10, 73628, 79670, 81180, 77404)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 3 + $i
    $values = $data[$i]
    $wsMem.Cells.Item($row, 1).Value = $values[0]
    $wsMem.Cells.Item($row, 2).Value = $values[1]
    $wsMem.Cells.Item($row, 3).Value = $values[2]
    $wsMem.Cells.Item($row, 4).Value = $values[3]
    $wsMem.Cells.Item($row, 5).Value = $values[4]
}

# Totals / averages row 13
$wsMem.Range("A13").Value = "Average"
$wsMem.Range("B13").Formula = "=AVERAGE(B3:B12)"
$wsMem.Range("C13").Formula = "=AVERAGE(C3:C12)"
$wsMem.Range("D13").Formula = "=AVERAGE(D3:D12)"
$wsMem.Range("E13").Formula = "=AVERAGE(E3:E12)"

# Column widths (bestFit, matching original formatting on similar sheets)
$wsMem.Columns.Item(3).ColumnWidth = 8.875
$wsMem.Columns.Item(5).ColumnWidth = 11.25

# Styling: copy styles from the equivalent cells on the JAVA SOCKETS RESULTS
# sheet so the new sheet visually matches its siblings.
$wsMem.Range("A1:E1").Style = $wsJava.Range("A1:E1").Style
$wsMem.Range("B2:E2").Style = $wsJava.Range("B2:E2").Style
$wsMem.Range("A3:A12").Style = $wsJava.Range("A3:A12").Style
$wsMem.Range("A13").Style = $wsJava.Range("A13").Style

# Selection on the new (now active) sheet
$wsMem.Range("G12").Select()
